$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.467.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.827.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5056"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3907"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07665"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.109"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.91%  "
$ws.Range("E12").Value = "  +3.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.283"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.579"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.823.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.30%  "
$ws.Range("E18").Value = "  +2.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06663"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.147"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.497.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.86%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.254"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  +2.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.034.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.396"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.125"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1081"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.675"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.663"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("E35").Value = "  +1.29%  "
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.963"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02320"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.135"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6245"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.16%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.180"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5899"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.720"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.979"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.44%  "
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06918"
$ws.Range("D51").Style = "Normal"
